# USERSTORIES.xlsx – remove the two "LOGOUT" user-story rows.
#
# The sheet lists FEATURES / SUB-FEATURES in a hierarchy of merged B-column
# group headers with C/D-column detail rows. Two "LOGOUT" entries need to be
# removed entirely:
#   1. Under ADMIN  -> LOGIN group: a standalone "LOGOUT" row pair (B19:B20).
#   2. Under USER    -> the trailing "LOGOUT" row pair at the end of that
#      section (originally rows 37:38, two rows above the EMPLOYEE header).
#
# Both deletions remove the rows entirely (shifting everything below up),
# matching the sharedStrings cleanup in the diff (LOGOUT / "  LOGOUT" /
# the long-padded "                         LOGOUT" strings all disappear).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) ADMIN section "LOGOUT" rows (B19 label, C20 value) – delete both rows.
$ws.Rows("19:20").Delete()

# 2) USER section "LOGOUT" rows. After the deletion above, these shifted up
#    from 37:38 to 35:36. The first of the pair (row 35, the "LOGOUT" label
#    in column B) becomes a blank spacer row that stays in place, while the
#    second (row 36, the "  LOGOUT" value in column C) is removed outright,
#    pulling the following EMPLOYEE section up by one more row.
$ws.Range("B35").ClearContents()
$ws.Rows("36:36").Delete()

# Restore the active selection shown in the saved view.
$ws.Range("C35").Select()
